$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.144.71"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "2.246.92"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'302.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.09%  "

$ws.Range("D6").Value = "'90.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.46%  "

$ws.Range("D7").Value = "'0.519"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.60%  "

$ws.Range("D10").Value = "'53.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.66%  "

$ws.Range("D11").Value = "'31.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.70%  "

$ws.Range("D12").Value = "'0.0793"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.31%  "

$ws.Range("E13").Value = "  +3.32%  "

$ws.Range("D14").Value = "'6.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.12%  "

$ws.Range("D15").Value = "2.591.20"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").Value = "'14.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").Value = "2.227.39"
$ws.Range("E17").Value = "  +1.29%  "

$ws.Range("D18").Value = "'0.748"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.44%  "

$ws.Range("D19").Value = "41.049.62"
$ws.Range("E19").Value = "  +3.05%  "

$ws.Range("D20").Value = "'11.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.71%  "

$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").Value = "'5.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.84%  "

$ws.Range("D23").Value = "'66.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.17%  "

$ws.Range("D24").Value = "'240.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("D25").Value = "'2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.48%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'1.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.50%  "

$ws.Range("D28").Value = "'23.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.47%  "

$ws.Range("D30").Value = "'9.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.84%  "

$ws.Range("D31").Value = "'157.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.80%  "

$ws.Range("D32").Value = "'33.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.25%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").Value = "'5.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.54%  "

$ws.Range("D35").Value = "'0.0731"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("D36").Value = "'3.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.72%  "

$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'16.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.79%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.115"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.93%  "

$ws.Range("D40").Value = "'0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.06%  "

$ws.Range("D41").Value = "'1.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.79%  "

$ws.Range("D42").Value = "'3.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.89%  "

$ws.Range("D43").Value = "2.065.87"
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("D44").Value = "'19.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.61%  "

$ws.Range("D45").Value = "'0.0276"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.76%  "

$ws.Range("E46").Value = "  +5.44%  "

$ws.Range("D47").Value = "'2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.28%  "

$ws.Range("E48").Value = "  -3.77%  "

$ws.Range("D49").Value = "2.459.81"
$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("D50").Value = "'1.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.44%  "

$ws.Range("D51").Value = "'1.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.79%  "
